# Updates cryptos list row values (coin name/link swaps + refreshed price/volume
# quotes) to match the latest GitHub Actions scrape, per commit
# "Updated cryptos list on Mon Jan 15 09:57:47 UTC 2024 with GitHub Actions".
#
# Column D prices that are plain decimals (e.g. "7.67") are written with a
# leading apostrophe so Excel keeps them as text instead of coercing them to
# numbers -- matching how the source data is stored (plain strings, some of
# which use "." as a thousands separator, e.g. "42.721.11").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.721.11"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3
$ws.Range("D3").Value = "2.534.20"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'316.46"
$ws.Range("E5").Value = "  +3.93%  "

# Row 6
$ws.Range("D6").Value = "'95.46"
$ws.Range("E6").Value = "  -2.92%  "

# Row 7
$ws.Range("E7").Value = "  +0.59%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  -1.16%  "

# Row 10
$ws.Range("D10").Value = "'36.34"
$ws.Range("E10").Value = "  -2.06%  "

# Row 11
$ws.Range("E11").Value = "  -1.58%  "

# Row 12
$ws.Range("D12").Value = "'7.67"
$ws.Range("E12").Value = "  -1.37%  "

# Row 13
$ws.Range("E13").Value = "  -1.87%  "

# Row 14
$ws.Range("D14").Value = "2.917.82"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("D15").Value = "'15.63"
$ws.Range("E15").Value = "  +3.51%  "

# Row 16
$ws.Range("D16").Value = "2.511.91"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("D17").Value = "'0.862"
$ws.Range("E17").Value = "  -1.12%  "

# Row 18
$ws.Range("D18").Value = "42.746.85"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  -5.59%  "

# Row 20
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0973"
$ws.Range("E21").Value = "  -2.04%  "

# Row 22
$ws.Range("D22").Value = "'71.43"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("D23").Value = "'254.76"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "'2.99"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25
$ws.Range("E25").Value = "  -1.44%  "

# Row 26
$ws.Range("D26").Value = "'27.70"
$ws.Range("E26").Value = "  -1.30%  "

# Row 27
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").Value = "'2.34"
$ws.Range("E28").Value = "  +12.25%  "

# Row 29
$ws.Range("D29").Value = "'39.91"
$ws.Range("E29").Value = "  +5.64%  "

# Row 30
$ws.Range("D30").Value = "'10.07"
$ws.Range("E30").Value = "  -1.92%  "

# Row 31
$ws.Range("D31").Value = "'5.92"
$ws.Range("E31").Value = "  -4.08%  "

# Row 32
$ws.Range("D32").Value = "'156.56"
$ws.Range("E32").Value = "  -1.28%  "

# Row 33
$ws.Range("D33").Value = "'20.04"
$ws.Range("E33").Value = "  +2.92%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'2.13"
$ws.Range("E34").Value = "  -0.78%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.35"
$ws.Range("E35").Value = "  +1.22%  "

# Row 36
$ws.Range("D36").Value = "'0.0789"

# Row 37
$ws.Range("E37").Value = "  -0.68%  "

# Row 38
$ws.Range("E38").Value = "  -3.22%  "

# Row 39
$ws.Range("D39").Value = "'24.90"
$ws.Range("E39").Value = "  -2.39%  "

# Row 40
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  +6.60%  "

# Row 42
$ws.Range("D42").Value = "'3.39"
$ws.Range("E42").Value = "  -1.47%  "

# Row 43
$ws.Range("D43").Value = "'3.86"

# Row 44
$ws.Range("E44").Value = "  -0.84%  "

# Row 45
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("D46").Value = "2.048.37"
$ws.Range("E46").Value = "  -2.44%  "

# Row 47
$ws.Range("D47").Value = "'86.04"
$ws.Range("E47").Value = "  -0.84%  "

# Row 48
$ws.Range("D48").Value = "'8.92"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'74.85"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.774.65"
$ws.Range("E50").Value = "  -0.80%  "

# Row 51
$ws.Range("E51").Value = "  -0.66%  "
